$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63. This shifts the existing rows 63-152
# down to 64-153, preserving all their values and formatting.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with a new price observation record
# (same Mercado/Region/Categoria/etc. as every other row in this sheet,
# with a new date and new price figures).
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112044
$ws.Range("G63").Value = "Perejil"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 180
$ws.Range("K63").Value = 4500
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = 4500
$ws.Range("N63").Value = "$/docena de atados (3 kilos)"
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 1500
$ws.Range("Q63").Value = 3
$ws.Range("R63").Value = "Hortaliza"
